$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 37 (pushes the old rows 37 "HOLY FAMILY MARONITE CHURCH"
# and 38 "SCHMITT MUSIC CTR" down to rows 38 and 39), matching the row's
# existing look (height/format) by giving it the same row height as its
# neighbours.
$ws.Rows.Item(37).Insert()
$ws.Rows.Item(37).RowHeight = 13.05

# Row 36 (MERWIN LIQUORS FALCON HEIGHTS) previously had a blank "Last Invoice
# Date" cell; fill it in with 12/05/2025 (serial 45996), reusing the date
# number format from the cell above so no new style gets created.
$ws.Range("D35").Copy()
$ws.Range("D36").PasteSpecial(-4122)
$ws.Range("D36").Value = 45996

# Populate the newly inserted row 37 with the new customer record.
$ws.Range("A37").Value = "HONEYCOMB SALON LLC"
$ws.Range("B37").Value = "Ballman, John W"
$ws.Range("C37").Value = "023"
$ws.Range("E37").Value = "0008385"
